# Updated cryptos list on Fri Nov  1 22:42:17 UTC 2024 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) (E) refresh, plus the
# Monero <-> Kaspa row-order swap (rows 36/37), exactly as described by the
# source diff.
#
# Note: several "Price" strings are plain decimal-looking text
# (e.g. "571.59") that Excel's COM layer would otherwise auto-coerce into a
# genuine number on assignment. A leading apostrophe forces those specific
# assignments to remain literal text, matching the workbook's original
# inline-string (text) cell type. Values that already contain two separators
# (e.g. "69.358.77") are not number-like and are assigned plainly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    # Force literal text even when $value looks like a plain number.
    $ws.Range($cell).Value = "'" + $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.358.77"
$ws.Range("E2").Value = "  -1.47%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.510.07"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "571.59"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6 - Solana
Set-TextValue "D6" "165.39"
$ws.Range("E6").Value = "  -1.78%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.03%  "

# Row 8 - XRP
Set-TextValue "D8" "0.512"
$ws.Range("E8").Value = "  +0.27%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.508.20"
$ws.Range("E9").Value = "  -0.40%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.36%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.40%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.357"
$ws.Range("E12").Value = "  +4.19%  "

# Row 13 - Toncoin
$ws.Range("E13").Value = "  +1.17%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.968.50"
$ws.Range("E14").Value = "  -0.05%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "69.150.10"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000175"
$ws.Range("E16").Value = "  -2.13%  "

# Row 17 - Avalanche
Set-TextValue "D17" "24.75"
$ws.Range("E17").Value = "  -0.72%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.526.56"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19 - Chainlink
Set-TextValue "D19" "11.26"
$ws.Range("E19").Value = "  -0.91%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -1.59%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "348.48"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -1.21%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +0.55%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.05%  "

# Row 25 - Litecoin
Set-TextValue "D25" "70.18"
$ws.Range("E25").Value = "  +1.73%  "

# Row 26 - NEARProtocol
$ws.Range("E26").Value = "  -4.35%  "

# Row 27 - Aptos
$ws.Range("E27").Value = "  -2.83%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.654.54"
$ws.Range("E28").Value = "  +0.35%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  +0.08%  "

# Row 30 - PEPE (0.0\u20830882 keeps the original subscript-3 character)
$ws.Range("D30").Value = "0.0₃0882"
$ws.Range("E30").Value = "  -3.14%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextValue "D31" "7.81"
$ws.Range("E31").Value = "  -0.21%  "

# Row 32 - Bittensor
Set-TextValue "D32" "460.09"
$ws.Range("E32").Value = "  -4.52%  "

# Row 33 - Fetch.AI
$ws.Range("E33").Value = "  -4.80%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  -2.16%  "

# Row 35 - FirstDigitalUSD
$ws.Range("E35").Value = "  +0.03%  "

# Rows 36/37 - Monero and Kaspa swap order (Kaspa now ranks above Monero)
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.117"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "157.27"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38 - WhiteBITCoin
$ws.Range("E38").Value = "  +0.88%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "18.41"
$ws.Range("E39").Value = "  -0.92%  "

# Row 41 - PolygonEcosystemToken
Set-TextValue "D41" "0.317"
$ws.Range("E41").Value = "  -0.95%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  -1.02%  "

# Row 43 - Stacks
Set-TextValue "D43" "1.60"
$ws.Range("E43").Value = "  -2.55%  "

# Row 44 - OKB
Set-TextValue "D44" "38.14"
$ws.Range("E44").Value = "  -0.39%  "

# Row 45 - ImmutableX
$ws.Range("E45").Value = "  -6.83%  "

# Row 46 - dogwifhat
Set-TextValue "D46" "2.22"
$ws.Range("E46").Value = "  -6.84%  "

# Row 47 - Aave
Set-TextValue "D47" "141.55"
$ws.Range("E47").Value = "  -0.95%  "

# Row 48 - Filecoin
Set-TextValue "D48" "3.45"
$ws.Range("E48").Value = "  -2.07%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  -1.81%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +0.18%  "

# Row 51 - Mantle
Set-TextValue "D51" "0.578"
$ws.Range("E51").Value = "  -3.39%  "
